# "Fecha de inicio de pago"
# Update the "Casos de Uso" task list: column T (row 13) now has 1 day of
# planned work logged, which ripples through the shared running-total
# formulas (U13, X13, AA13, ... AY13, AZ13, BA13) across the rest of the
# row via the existing shared formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Log 1 day of work in the T13 cell (previously empty); downstream shared
# formulas (U13=R13-T13, X13=U13-W13, ... BA13=G13-AZ13) recalculate
# automatically.
$ws.Range("T13").Value = 1

# Reflect the author's final cell selection on that row.
$ws.Range("W13").Select()
